# Commit: Thu, Jun 04, 2020  6:05:16 PM
#
# The canonical diff swaps the contents of ppt/theme/theme1.xml (used by
# the slide master / "Integral" design) and ppt/theme/theme2.xml (used by
# the notes master / "Office Theme" design) - i.e. the slide master's
# design switches from the "Integral" (Red Violet) palette to the plain
# "Office" palette, while the notes master picks up the "Integral"
# palette.
#
# Only the <a:clrScheme> (12 colour slots) differs between the two theme
# parts - the font scheme and format scheme (fills/lines/effects) are
# byte-identical in both, so re-pointing the 12 theme colours on the
# slide master's Theme reproduces the target theme1.xml.
#
# PowerPoint's object model doesn't expose a scriptable "apply/replace
# theme" verb, so this goes through the same automation surface real
# VBA macros use for recolouring a theme:
#   ActivePresentation.SlideMaster.Theme.ThemeColorScheme.Colors(i).RGB = ...
# with i in msoThemeColorDark1..msoThemeColorFollowedHyperlink order
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

function Set-ThemeColor($index, $red, $green, $blue) {
    $colorScheme.Colors($index).RGB = $red + ($green * 256) + ($blue * 65536)
}

# New ("Office Theme") palette - matches ppt/theme/theme2.xml's <a:clrScheme>
Set-ThemeColor 1  0x00 0x00 0x00   # dk1
Set-ThemeColor 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor 12 0x95 0x4F 0x72   # folHlink
